# Update the top-level-description (column B) text for several rows in the
# "Specification" worksheet to the newer, shorter descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Specification")

$updates = @{
    21 = "Name and contact information if an agent is being used."
    25 = "Name and contact information if an agent is being used."
    33 = "Telephone number and email address of the applicant."
    37 = "Name and contact information for the parties making the application."
    43 = "Details of any conflict of interest that may exist between the applicant and planning authority."
    46 = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
    47 = "Signed and dated verification of the application's accuracy."
    50 = "Names and contact details for all parties with an interest in the proposed develpoment."
    68 = "Details of pre-application advice received from the planning authority"
    73 = "Where the proposed development will be built."
    82 = "Information to help the planning authority arrange a site visit"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
